$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "commands" sheet (sheet1): restructure to add a "Type" column and a
#    new "pwd" row, and append the git command table that used to live on
#    the "git" sheet.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("commands")

# Insert a new row above the old row 7 ("cacls/chmod") to make room for the
# new "pwd" entry; everything below shifts down by one row.
$ws1.Rows.Item(7).Insert()

# Insert a new column before the old "Description" column (D) so it becomes
# column E, leaving room for the new "Type" column in D.
$ws1.Columns.Item(4).Insert()

# New column widths (A-D); column E keeps the width that travelled with the
# shifted "Description" column.
$ws1.Columns.Item(1).ColumnWidth = 28.42
$ws1.Columns.Item(2).ColumnWidth = 11.42
$ws1.Columns.Item(3).ColumnWidth = 9.59
$ws1.Columns.Item(4).ColumnWidth = 11.59

# Header row.
$ws1.Range("D1").Value = "Type"

# New "pwd" row.
$ws1.Range("A7").Value = "pwd"
$ws1.Range("E7").Value = "print working directory"

# "Type" tag for the netstat row (now row 11).
$ws1.Range("D11").Value = "network"

# Git command table (previously on the "git" sheet), appended below the
# existing data.
$ws1.Range("A16").Value = "git clone url"
$ws1.Range("D16").Value = "git"
$ws1.Range("E16").Value = "clones a github repository"

$ws1.Range("A17").Value = "git remote -v"
$ws1.Range("D17").Value = "git"
$ws1.Range("E17").Value = "displays the path to the remote origin"

$ws1.Range("A18").Value = "git status"
$ws1.Range("D18").Value = "git"
$ws1.Range("E18").Value = "shows the status of the directory"

$ws1.Range("A19").Value = "git reset HEAD myFile.txt"
$ws1.Range("D19").Value = "git"
$ws1.Range("E19").Value = "unstages the file "

$ws1.Range("A20").Value = "git checkout HEAD myFile.txt"
$ws1.Range("D20").Value = "git"
$ws1.Range("E20").Value = "gets the latest committed version of the file"

$ws1.Range("A21").Value = "git add ."
$ws1.Range("D21").Value = "git"
$ws1.Range("E21").Value = "stages all modified files "

$ws1.Range("A22").Value = "git commit -a"
$ws1.Range("D22").Value = "git"
$ws1.Range("E22").Value = "commits modified files, stages added and removed files"

$ws1.Range("E7").Select() | Out-Null

# ---------------------------------------------------------------------------
# 2. "git" sheet (sheet2): its table moved onto the "commands" sheet, so wipe
#    its contents, resize column A and hide the sheet.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("git")
$ws2.Range("A1:B7").ClearContents()
$ws2.Columns.Item(1).ColumnWidth = 68.59
$ws2.Range("B1:B7").Select() | Out-Null
$ws2.Visible = $false

# ---------------------------------------------------------------------------
# 3. "Sheet3": no longer used, remove it entirely.
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Sheet3")
$ws3.Delete() | Out-Null

$ws1.Activate()
